<#
.SYNOPSIS
    Renames the workbook's single worksheet to a descriptive name and
    restores the author's last-used cell selection.

.DESCRIPTION
    The sheet was generically named "Sheet1"; rename it to
    "cover_type_keys" to reflect the lookup-table data it holds
    (cover-type keys used by the classifyTypeCover configuration).
    Also moves the active selection to I31 (the cell that was selected
    when the workbook was last saved).
#>

try {
    $wb = $excel.ActiveWorkbook
    if (-not $wb) {
        throw "No active workbook found."
    }

    $ws = $wb.ActiveSheet
    if (-not $ws) {
        throw "No active worksheet found on the workbook."
    }

    # Give the sheet a meaningful name instead of the generic default.
    $ws.Name = "cover_type_keys"

    # Restore the last-used selection on the sheet.
    $ws.Range("I31").Select()
}
catch {
    Write-Error "edit.ps1 failed while updating '$($wb.Name)': $_"
    throw
}
